# Add data for 2022-12-14: update the "as of" date from 12-05 to 12-06,
# and update December (row 13) and Total (row 14) carjacking counts for
# the "2022 (through ...)" column (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and the header label that both encode the "as of" date.
$ws.Name = "Through 2022-12-06"
$ws.Range("I1").Value = "2022 (through 12-06)"

# Update December 2022 carjacking count.
$ws.Range("I13").Value = 24

# Update the yearly Total for the 2022 column.
$ws.Range("I14").Value = 1539
